# This script replicates a Power Query "Refresh" of the 쿼리1 (Query1) table
# on the first worksheet: updated ranking data (column B/C) and a new
# refresh timestamp (column D), plus Excel's explicit re-application of the
# "General" number format to the text column during the refresh, and an
# updated active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New refresh timestamp applied to every data row (D2:D11)
$newTimestamp = 46025.584667777781

# New ranking data: row -> (name, value)
$names  = @("태영", "하루묭", "빵지니", "으냉이", "임밍지", "한쪼니", "우리밍", "윤하랑", "히요코", "하랑e")
$values = @(110420, 51985, 48927, 43645, 31707, 28546, 26528, 17690, 2154, $null)

for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i

    $ws.Cells.Item($row, 2).Value = $names[$i]

    if ($null -eq $values[$i]) {
        $ws.Cells.Item($row, 3).Value = ""
    } else {
        $ws.Cells.Item($row, 3).Value = $values[$i]
    }

    $ws.Cells.Item($row, 4).Value = $newTimestamp
}

# Excel's query-table refresh re-applies an explicit "General" number format
# to the refreshed text column (column B), producing a new explicit style.
$ws.Range("B2:B11").NumberFormat = "General"

# Update the saved selection to match the post-refresh cursor position.
$ws.Range("G17").Select()
